$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = "[0.1912016322809631, 12.821937688449307]"
$ws.Range("N2").Value = [double]"0.04372556087122881"
$ws.Range("O2").Value = [double]"0.04372556087122881"
$ws.Range("Q2").Value = "[-2.7925268031909276, -0.1006315965113842]"
$ws.Range("R2").Value = [double]"0.03575428185013729"
$ws.Range("S2").Value = [double]"0.03575428185013729"
$ws.Range("U2").Value = "[4.97034028870712, 12.621294173096835]"
$ws.Range("V2").Value = [double]"3.10863189627586e-05"
$ws.Range("W2").Value = [double]"3.10863189627586e-05"
$ws.Range("Y2").Value = [double]"0.4143343343343417"
$ws.Range("Z2").Value = [double]"11.49777777777805"

# Row 3
$ws.Range("M3").Value = "[-1.0993139906444487, 13.860363846513692]"
$ws.Range("N3").Value = [double]"0.09265477698908731"
$ws.Range("O3").Value = [double]"0.09265477698908731"
$ws.Range("U3").Value = "[5.489874854514014, 13.832642375173808]"
$ws.Range("V3").Value = [double]"2.78291921482765e-05"
$ws.Range("W3").Value = [double]"2.78291921482765e-05"

# Row 4
$ws.Range("M4").Value = "[-1.9189362516150794, 14.72835311345948]"
$ws.Range("N4").Value = [double]"0.1282015086036283"
$ws.Range("O4").Value = [double]"0.1282015086036283"
$ws.Range("Q4").Value = "[-5.522158858562238, 0.49057903299299976]"
$ws.Range("R4").Value = [double]"0.09882671305278001"
$ws.Range("S4").Value = [double]"0.09882671305278001"
$ws.Range("U4").Value = "[5.259193296356727, 13.683796185967743]"
$ws.Range("V4").Value = [double]"4.335989771764304e-05"
$ws.Range("W4").Value = [double]"4.335989771764304e-05"
$ws.Range("Y4").Value = [double]"-2.019879879879923"
$ws.Range("Z4").Value = [double]"22.73659659659712"

# Row 5
$ws.Range("M5").Value = "[-0.6925043185537447, 13.732060425444757]"
$ws.Range("N5").Value = [double]"0.07530376671966543"
$ws.Range("O5").Value = [double]"0.07530376671966543"
$ws.Range("Q5").Value = "[1.603816069400195, 4.597606065613891]"
$ws.Range("R5").Value = [double]"0.0001359863332710542"
$ws.Range("S5").Value = [double]"0.0001359863332710542"
$ws.Range("U5").Value = "[5.231407146176435, 13.691033822150851]"
$ws.Range("V5").Value = [double]"4.68150624550745e-05"
$ws.Range("W5").Value = [double]"4.68150624550745e-05"
$ws.Range("Y5").Value = [double]"6.940100100100263"
$ws.Range("Z5").Value = [double]"19.26654654654699"

# Row 6
$ws.Range("M6").Value = "[-0.8069782563330623, 13.693337068055216]"
$ws.Range("N6").Value = [double]"0.0801963029959063"
$ws.Range("O6").Value = [double]"0.0801963029959063"
$ws.Range("Q6").Value = "[1.4277107755052714, 4.647921863869584]"
$ws.Range("R6").Value = [double]"0.0004319077789900483"
$ws.Range("S6").Value = [double]"0.0004319077789900483"
$ws.Range("U6").Value = "[5.259680032142127, 13.695479305548597]"
$ws.Range("V6").Value = [double]"4.379903471041224e-05"
$ws.Range("W6").Value = [double]"4.379903471041224e-05"
$ws.Range("Y6").Value = [double]"6.732932932933092"
$ws.Range("Z6").Value = [double]"19.9916316316321"

# Row 7
$ws.Range("M7").Value = "[-1.4760338678708855, 14.221486347130059]"
$ws.Range("N7").Value = [double]"0.108956252782286"
$ws.Range("O7").Value = [double]"0.108956252782286"
$ws.Range("Q7").Value = "[-0.767315923399309, 5.484422009870469]"
$ws.Range("R7").Value = [double]"0.1355831019938456"
$ws.Range("S7").Value = [double]"0.1355831019938456"
$ws.Range("U7").Value = "[5.511785641543515, 13.928659396309644]"
$ws.Range("V7").Value = [double]"2.902433830431761e-05"
$ws.Range("W7").Value = [double]"2.902433830431761e-05"
$ws.Range("Y7").Value = [double]"3.288778778778855"
$ws.Range("Z7").Value = [double]"29.02929929929998"

# Row 8
$ws.Range("M8").Value = "[-0.3951469098465026, 13.812938491973817]"
$ws.Range("N8").Value = [double]"0.06357395869666771"
$ws.Range("O8").Value = [double]"0.06357395869666771"
$ws.Range("Q8").Value = "[0.47171060864711567, 3.4403427057329674]"
$ws.Range("R8").Value = [double]"0.01095033289426883"
$ws.Range("S8").Value = [double]"0.01095033289426883"
$ws.Range("U8").Value = "[5.092489021505877, 12.717454614548746]"
$ws.Range("V8").Value = [double]"2.444259258016501e-05"
$ws.Range("W8").Value = [double]"2.444259258016501e-05"
$ws.Range("Y8").Value = [double]"11.70494494494522"
$ws.Range("Z8").Value = [double]"23.92780780780837"

# Row 9
$ws.Range("M9").Value = "[-0.3741280593200056, 12.7217586624634]"
$ws.Range("N9").Value = [double]"0.06398367038865316"
$ws.Range("O9").Value = [double]"0.06398367038865316"
$ws.Range("Q9").Value = "[0.42139481039142357, 3.4906585039886595]"
$ws.Range("R9").Value = [double]"0.0136487085209076"
$ws.Range("S9").Value = [double]"0.0136487085209076"
$ws.Range("U9").Value = "[3.4110932778948477, 10.372734283033935]"
$ws.Range("V9").Value = [double]"0.0002422441177527546"
$ws.Range("W9").Value = [double]"0.0002422441177527546"
$ws.Range("Y9").Value = [double]"11.49777777777805"
$ws.Range("Z9").Value = [double]"24.13497497497554"

# Row 10
$ws.Range("M10").Value = "[-0.7199709953569808, 14.222150354962857]"
$ws.Range("N10").Value = [double]"0.07541308875262165"
$ws.Range("O10").Value = [double]"0.07541308875262165"
$ws.Range("Q10").Value = "[-0.19497371824080822, 2.8994478744842738]"
$ws.Range("R10").Value = [double]"0.0851537407081111"
$ws.Range("S10").Value = [double]"0.0851537407081111"
$ws.Range("U10").Value = "[4.610227996328369, 12.81486114793205]"
$ws.Range("V10").Value = [double]"9.72909759797691e-05"
$ws.Range("W10").Value = [double]"9.72909759797691e-05"
$ws.Range("Y10").Value = [double]"13.93199199199232"
$ws.Range("Z10").Value = [double]"26.6727727727734"

# Row 11
$ws.Range("M11").Value = "[1.3805700443853723, 16.2008889326949]"
$ws.Range("N11").Value = [double]"0.02113010238668411"
$ws.Range("O11").Value = [double]"0.02113010238668411"
$ws.Range("Q11").Value = "[-0.20755266780472947, 1.4025528763774249]"
$ws.Range("R11").Value = [double]"0.1419357197691762"
$ws.Range("S11").Value = [double]"0.1419357197691762"
$ws.Range("U11").Value = "[6.008279532594768, 13.91319160224335]"
$ws.Range("V11").Value = [double]"7.143718056434878e-06"
$ws.Range("W11").Value = [double]"7.143718056434878e-06"
$ws.Range("Y11").Value = [double]"20.09521521521568"
$ws.Range("Z11").Value = [double]"26.72456456456518"

# Row 12
$ws.Range("M12").Value = "[-2.476004497284068, 12.921516058464507]"
$ws.Range("N12").Value = [double]"0.1786195058833675"
$ws.Range("O12").Value = [double]"0.1786195058833675"
$ws.Range("Q12").Value = "[-2.7736583788450435, 3.4906585039886586]"
$ws.Range("R12").Value = [double]"0.8187248789934645"
$ws.Range("S12").Value = [double]"0.8187248789934645"
$ws.Range("U12").Value = "[4.445485625523846, 12.630090952619026]"
$ws.Range("V12").Value = [double]"0.0001236875072767596"
$ws.Range("W12").Value = [double]"0.0001236875072767596"
$ws.Range("Y12").Value = [double]"9.973333333333366"
$ws.Range("Z12").Value = [double]"32.34594594594605"

# Row 13
$ws.Range("M13").Value = "[-0.16612312703337828, 13.442026483610556]"
$ws.Range("N13").Value = [double]"0.05561431441904729"
$ws.Range("O13").Value = [double]"0.05561431441904729"
$ws.Range("Q13").Value = "[-1.1446844103170015, 1.9120003337163105]"
$ws.Range("R13").Value = [double]"0.6156066229928951"
$ws.Range("S13").Value = [double]"0.6156066229928951"
$ws.Range("U13").Value = "[5.3148709064602215, 13.239151599529293]"
$ws.Range("V13").Value = [double]"2.354303029150628e-05"
$ws.Range("W13").Value = [double]"2.354303029150628e-05"
$ws.Range("Y13").Value = [double]"15.61141141141146"
$ws.Range("Z13").Value = [double]"26.52816816816825"

# Row 14
$ws.Range("M14").Value = "[-0.47941862656526, 16.355828623208968]"
$ws.Range("N14").Value = [double]"0.0639341405519267"
$ws.Range("O14").Value = [double]"0.0639341405519267"
$ws.Range("Q14").Value = "[-2.2705003962881176, 3.9183427891620437]"
$ws.Range("R14").Value = [double]"0.5944094201815537"
$ws.Range("S14").Value = [double]"0.5944094201815537"
$ws.Range("U14").Value = "[6.564110989370478, 15.280638626870838]"
$ws.Range("V14").Value = [double]"7.849340683874217e-06"
$ws.Range("W14").Value = [double]"7.849340683874217e-06"
$ws.Range("Y14").Value = [double]"8.445885885885914"
$ws.Range("Z14").Value = [double]"30.54894894894904"

# Row 15
$ws.Range("M15").Value = "[-0.34596936915776233, 14.643612281458145]"
$ws.Range("N15").Value = [double]"0.06106379566448394"
$ws.Range("O15").Value = [double]"0.06106379566448394"
$ws.Range("Q15").Value = "[-0.06918422260157708, 2.899447874484273]"
$ws.Range("R15").Value = [double]"0.06118093329988983"
$ws.Range("S15").Value = [double]"0.06118093329988983"
$ws.Range("U15").Value = "[4.979105169123802, 13.309521704690656]"
$ws.Range("V15").Value = [double]"6.128749666189925e-05"
$ws.Range("W15").Value = [double]"6.128749666189925e-05"
$ws.Range("Y15").Value = [double]"12.08480480480484"
$ws.Range("Z15").Value = [double]"22.68708708708716"
